$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy formatting (styles) from row 841 down through new rows 842-857
$ws.Range("A841:I841").Copy()
$ws.Range("A842:I857").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Re-apply the blank-G-cell style (as used on e.g. G836) to rows where G ends up empty
$ws.Range("G836").Copy()
$ws.Range("G844").PasteSpecial(-4122)
$ws.Range("G836").Copy()
$ws.Range("G848").PasteSpecial(-4122)
$ws.Range("G836").Copy()
$ws.Range("G850").PasteSpecial(-4122)
$ws.Range("G836").Copy()
$ws.Range("G851").PasteSpecial(-4122)
$ws.Range("G836").Copy()
$ws.Range("G853").PasteSpecial(-4122)
$ws.Range("G836").Copy()
$ws.Range("G856").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 842
$ws.Range("A842").Value = 46070
$ws.Range("B842").Value = "Kamal Bafounta"
$ws.Range("C842").Value = 70
$ws.Range("D842").Value = 7
$ws.Range("E842").Value = 3
$ws.Range("F842").Value = 1
$ws.Range("G842").Value = "Genou/cheville"
$ws.Range("H842").Value = 7

# Row 843
$ws.Range("A843").Value = 46070
$ws.Range("B843").Value = "Yoan Zouma"
$ws.Range("C843").Value = 70
$ws.Range("D843").Value = 5
$ws.Range("E843").Value = 7
$ws.Range("F843").Value = 3
$ws.Range("G843").Value = "Cheville droite "
$ws.Range("H843").Value = 5

# Row 844
$ws.Range("A844").Value = 46070
$ws.Range("B844").Value = "Omar Benyounes"
$ws.Range("C844").Value = 70
$ws.Range("D844").Value = 7
$ws.Range("E844").Value = 5
$ws.Range("F844").Value = 0
$ws.Range("H844").Value = 8

# Row 845
$ws.Range("A845").Value = 46070
$ws.Range("B845").Value = "Romain Thunet"
$ws.Range("C845").Value = 70
$ws.Range("D845").Value = 7
$ws.Range("E845").Value = 6
$ws.Range("F845").Value = 3
$ws.Range("G845").Value = "Synthétique "
$ws.Range("H845").Value = 7

# Row 846
$ws.Range("A846").Value = 46070
$ws.Range("B846").Value = "Naim Ighbane"
$ws.Range("C846").Value = 70
$ws.Range("D846").Value = 7
$ws.Range("E846").Value = 6
$ws.Range("F846").Value = 6
$ws.Range("G846").Value = "Genou"
$ws.Range("H846").Value = 7

# Row 847
$ws.Range("A847").Value = 46070
$ws.Range("B847").Value = "Yoann Martelat"
$ws.Range("C847").Value = 70
$ws.Range("D847").Value = 6
$ws.Range("E847").Value = 4
$ws.Range("F847").Value = 5
$ws.Range("G847").Value = "Genou"
$ws.Range("H847").Value = 6

# Row 848
$ws.Range("A848").Value = 46070
$ws.Range("B848").Value = "Malik Boussaid"
$ws.Range("C848").Value = 70
$ws.Range("D848").Value = 5
$ws.Range("E848").Value = 3
$ws.Range("F848").Value = 0
$ws.Range("H848").Value = 10

# Row 849
$ws.Range("A849").Value = 46070
$ws.Range("B849").Value = "Mehdi Boussaid"
$ws.Range("C849").Value = 70
$ws.Range("D849").Value = 6
$ws.Range("E849").Value = 7
$ws.Range("F849").Value = 3
$ws.Range("G849").Value = "Adducteur "
$ws.Range("H849").Value = 7

# Row 850
$ws.Range("A850").Value = 46070
$ws.Range("B850").Value = "Levy Ndoutoume"
$ws.Range("C850").Value = 70
$ws.Range("D850").Value = 7
$ws.Range("E850").Value = 7
$ws.Range("F850").Value = 0
$ws.Range("H850").Value = 6

# Row 851
$ws.Range("A851").Value = 46070
$ws.Range("B851").Value = "Ilan Ihaddadene"
$ws.Range("C851").Value = 70
$ws.Range("D851").Value = 7
$ws.Range("E851").Value = 7
$ws.Range("F851").Value = 0
$ws.Range("H851").Value = 6

# Row 852
$ws.Range("A852").Value = 46070
$ws.Range("B852").Value = "Karahali Souaré"
$ws.Range("C852").Value = 70
$ws.Range("D852").Value = 8
$ws.Range("E852").Value = 7
$ws.Range("F852").Value = 6
$ws.Range("G852").Value = "Cheville "
$ws.Range("H852").Value = 1

# Row 853
$ws.Range("A853").Value = 46070
$ws.Range("B853").Value = "Theo Owono"
$ws.Range("C853").Value = 70
$ws.Range("D853").Value = 8
$ws.Range("E853").Value = 7
$ws.Range("F853").Value = 0
$ws.Range("H853").Value = 6

# Row 854
$ws.Range("A854").Value = 46070
$ws.Range("B854").Value = "Nathanael Beta"
$ws.Range("C854").Value = 70
$ws.Range("D854").Value = 5
$ws.Range("E854").Value = 5
$ws.Range("F854").Value = 1
$ws.Range("G854").Value = "Courbatures "
$ws.Range("H854").Value = 6

# Row 855
$ws.Range("A855").Value = 46070
$ws.Range("B855").Value = "Sofiane Belle"
$ws.Range("C855").Value = 70
$ws.Range("D855").Value = 4
$ws.Range("E855").Value = 5
$ws.Range("F855").Value = 4
$ws.Range("G855").Value = "Ventre"
$ws.Range("H855").Value = 6

# Row 856
$ws.Range("A856").Value = 46070
$ws.Range("B856").Value = "Mattheo Haon"
$ws.Range("C856").Value = 70
$ws.Range("D856").Value = 8
$ws.Range("E856").Value = 5
$ws.Range("F856").Value = 0
$ws.Range("H856").Value = 7

# Row 857
$ws.Range("A857").Value = 46070
$ws.Range("B857").Value = "Jeremie Laurent"
$ws.Range("C857").Value = 70
$ws.Range("D857").Value = 8
$ws.Range("E857").Value = 7
$ws.Range("F857").Value = 1
$ws.Range("G857").Value = "Courbaturé "
$ws.Range("H857").Value = 7

# Extend the C*D formula down through the new rows (creates a shared formula group)
$ws.Range("I842:I857").Formula = "=C842*D842"

# Update the view: selection on K852 (matches target sheetView state)
[void]$ws.Range("K852").Select()
